$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.636.96'
$ws.Range('E2').Value = '  +2.59%  '
$ws.Range('D3').Value = '2.432.65'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '565.16'
$ws.Range('E5').Value = '  +2.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.92'
$ws.Range('E6').Value = '  +5.77%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.514'
$ws.Range('E8').Value = '  +2.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.171'
$ws.Range('E9').Value = '  +8.44%  '
$ws.Range('D10').Value = '2.431.38'
$ws.Range('E10').Value = '  +1.39%  '
$ws.Range('E11').Value = '  -1.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.334'
$ws.Range('E12').Value = '  +2.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.70'
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('E14').Value = '  +6.17%  '
$ws.Range('D15').Value = '69.640.74'
$ws.Range('E15').Value = '  +2.73%  '
$ws.Range('D16').Value = '2.881.03'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.05'
$ws.Range('E17').Value = '  +5.68%  '
$ws.Range('D18').Value = '2.442.25'
$ws.Range('E18').Value = '  +1.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.88'
$ws.Range('E19').Value = '  +5.86%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.20'
$ws.Range('E20').Value = '  +6.07%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '344.06'
$ws.Range('E21').Value = '  +4.29%  '
$ws.Range('E22').Value = '  +3.07%  '
$ws.Range('E23').Value = '  +8.42%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.96'
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.83'
$ws.Range('E26').Value = '  +5.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.50'
$ws.Range('E27').Value = '  +5.72%  '
$ws.Range('D28').Value = '2.561.64'
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('D30').Value = '0.0₃0857'
$ws.Range('E30').Value = '  +8.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.38'
$ws.Range('E31').Value = '  +4.95%  '
$ws.Range('E32').Value = '  +10.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '457.33'
$ws.Range('E33').Value = '  +8.99%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('E35').Value = '  +2.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.00'
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range('B37').Value = 'WhiteBITCoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.13'
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.112'
$ws.Range('E38').Value = '  +7.41%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.27'
$ws.Range('E40').Value = '  +3.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.304'
$ws.Range('E41').Value = '  +4.27%  '
$ws.Range('E42').Value = '  +4.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.52'
$ws.Range('E43').Value = '  +4.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.07'
$ws.Range('E44').Value = '  +1.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.08'
$ws.Range('E45').Value = '  +3.41%  '
$ws.Range('E46').Value = '  +9.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '136.03'
$ws.Range('E47').Value = '  +5.62%  '
$ws.Range('E48').Value = '  +3.78%  '
$ws.Range('E50').Value = '  +3.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.564'
$ws.Range('E51').Value = '  +2.09%  '
